# Fix lookup table entries ("Some letters were not correctly matching.").
# Each row maps a letter/glyph to Top/Middle/Bottom (+ flip variants) match
# flags; the corrections below re-point several rows at the right glyph.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Letter "d" (row 5): Top_hflip should be set.
$ws.Range("E5").Value = 1

# Letter "j" (row 11): Top_hflip should be set.
$ws.Range("E11").Value = 1

# Letter "C" (row 30): Top match should point at glyph 31, not 33.
$ws.Range("B30").Value = 31

# Letter "G" (row 34): Top match should point at glyph 31, not 4.
$ws.Range("B34").Value = 31

# Letter "S" (row 46): Top match should point at glyph 31, not 33.
$ws.Range("B46").Value = 31

# Letter "T" (row 47): Top match should point at glyph 33, not 31.
$ws.Range("B47").Value = 33

# Letter "Z" (row 53): Top match should point at glyph 33, not 31.
$ws.Range("B53").Value = 33

# "(" (row 66): Top match should point at glyph 3, not 31.
$ws.Range("B66").Value = 3

# ")" (row 67): Top match should point at glyph 3, not 31.
$ws.Range("B67").Value = 3

# "?" (row 70): Middle match should point at glyph 76, not 78.
$ws.Range("C70").Value = 76

# Leave the view scrolled/selected where the author last left it.
$ws.Range("E51").Select()
